$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Third Iteration")

# Row 15: Dead_Temp_sensor_alert
$ws.Range("A15").Value = "Dead_Temp_sensor_alert"
$ws.Range("B15").Value = "Temperature"
$ws.Range("C15").Value = "*"
$ws.Range("D15").Value = "deg-F"
$ws.Range("E15").Value = "<alias> has not reported recently."
$ws.Range("F15").Value = "WATT"
$ws.Range("G15").Value = "Watt"
$ws.Range("H15").Value = "CEVAC_WATT_TEMP_LATEST"
$ws.Range("I15").Value = "UTCDateTime"
$ws.Range("J15").Value = "UTCDateTime"
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = "*"
$ws.Range("N15").Value = "<"
$ws.Range("O15").Value = "<now> - 1 hr"
$ws.Range("P15").Value = "alert"
$ws.Range("Q15").Value = "subtraction and alert value are in hours"

# Row 16: Dead_Power_sensor_alert
$ws.Range("A16").Value = "Dead_Power_sensor_alert"
$ws.Range("B16").Value = "Power"
$ws.Range("C16").Value = "*"
$ws.Range("D16").Value = "kW"
$ws.Range("E16").Value = "<alias> has not reported recently."
$ws.Range("F16").Value = "WATT"
$ws.Range("G16").Value = "Watt"
$ws.Range("H16").Value = "CEVAC_WATT_POWER_LATEST"
$ws.Range("I16").Value = "UTCDateTime"
$ws.Range("J16").Value = "UTCDateTime"
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = "*"
$ws.Range("N16").Value = "<"
$ws.Range("O16").Value = "<now> - 1 hr"
$ws.Range("P16").Value = "alert"

# Row 17: Dead_IAQ_sensor_alert
$ws.Range("A17").Value = "Dead_IAQ_sensor_alert"
$ws.Range("B17").Value = "CO2"
$ws.Range("C17").Value = "*"
$ws.Range("D17").Value = "ppm"
$ws.Range("E17").Value = "<alias> has not reported recently."
$ws.Range("F17").Value = "WATT"
$ws.Range("G17").Value = "Watt"
$ws.Range("H17").Value = "CEVAC_WATT_IAQ_LATEST"
$ws.Range("I17").Value = "UTCDateTime"
$ws.Range("J17").Value = "UTCDateTime"
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = "*"
$ws.Range("N17").Value = "<"
$ws.Range("O17").Value = "<now> - 1 hr"
$ws.Range("P17").Value = "alert"

# Update selection to match post-edit state
$ws.Range("A18").Select()

Write-Host "done"
